# Append two new trading-log rows (92 and 93) to Sheet1, matching the
# 2025-10-20 12:44:5x UTC BTC trade that was logged upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 92: TRADING_ATTEMPT
$ws.Cells.Item(92, 1).Value = "2025-10-20T12:44:52.001924"
$ws.Cells.Item(92, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(92, 3).Value = "BTC"
$ws.Cells.Item(92, 4).Value = "UNKNOWN"
$ws.Cells.Item(92, 5).Value = 108621.134857428
$ws.Cells.Item(92, 11).Value = "ATTEMPT"
$ws.Cells.Item(92, 12).Value = "Attempting trade 1/1"

# Row 93: POSITION_OPENED
$ws.Cells.Item(93, 1).Value = "2025-10-20T12:44:53.326525"
$ws.Cells.Item(93, 2).Value = "POSITION_OPENED"
$ws.Cells.Item(93, 3).Value = "BTC"
$ws.Cells.Item(93, 4).Value = "UNKNOWN"
$ws.Cells.Item(93, 5).Value = 108621.134857428
$ws.Cells.Item(93, 6).Value = 3600
$ws.Cells.Item(93, 7).Value = 40
$ws.Cells.Item(93, 8).Value = 0.08208394928824908
$ws.Cells.Item(93, 11).Value = "SUCCESS"
